$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add Pruebas Realizadas value for the existing last row (row 44, date 43963 / 2020-05-12)
$ws.Range("B44").Value = 529

# Add a new row (row 45) with the next day's date and its "Pruebas Realizadas" value
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A45").Value = 43964
$ws.Range("B45").Value = 803

# Extend the Excel table "Condicion_Pacientes" so it covers the newly added row
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F45"))

# Update the active selection to match the post-edit workbook state
$ws.Range("C45").Select()
